$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 11 (20190710): fill in the checklist results for columns B..O
$ws.Range("B11").Value = "×"
$ws.Range("C11").Value = "√"
$ws.Range("D11").Value = "√"
$ws.Range("E11").Value = "×"
$ws.Range("F11").Value = "×"
$ws.Range("G11").Value = "√"
$ws.Range("H11").Value = "×"
$ws.Range("I11").Value = "√"
$ws.Range("J11").Value = "√"
$ws.Range("K11").Value = "√"
$ws.Range("L11").Value = "×"
$ws.Range("M11").Value = "√"
$ws.Range("N11").Value = "×"
$ws.Range("O11").Value = "×"

# Row 12 (20190711): new row, initialized with mostly "√" (D12 and H12 left blank)
$ws.Range("A12").Value = 20190711
$ws.Range("B12").Value = "√"
$ws.Range("C12").Value = "√"
$ws.Range("E12").Value = "√"
$ws.Range("F12").Value = "√"
$ws.Range("G12").Value = "√"
$ws.Range("I12").Value = "√"
$ws.Range("J12").Value = "√"
$ws.Range("K12").Value = "√"
$ws.Range("L12").Value = "√"
$ws.Range("M12").Value = "√"
$ws.Range("N12").Value = "√"
$ws.Range("O12").Value = "√"

$ws.Range("O12").Select()
